# Adds generic code for the ship and bullets (spawn, move, remove, etc.)
# to the "Asteroids" section of the Subroutines sheet, reordering the
# existing asteroid-related rows and inserting 8 new rows for the new
# ship/bullet subroutines, then updates the sheet view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subroutines")

# --- Insert 8 fresh rows right after "initializeShip" (row 34), before
#     the existing "drawAsteroids" row (row 35). This shifts every row
#     at/after 35 down by 8 and carries the column B fill formatting
#     (style s="8") down into the newly inserted rows automatically.
$ws.Rows("35:42").Insert()

# --- Mark A33 like the other blank-but-formatted A cells in this sheet
#     (copy formatting from A10, which already carries that same style).
$ws.Range("A10").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the 8 new "ship" / "bullet" subroutine rows (35-42).
$ws.Range("B35").Value = "moveShip"
$ws.Range("B36").Value = "drawShip"
$ws.Range("B37").Value = "createBullet"
$ws.Range("B38").Value = "moveBullets"
$ws.Range("B39").Value = "drawBullets"
$ws.Range("B40").Value = "flashBullets"
$ws.Range("B41").Value = "removeBullet"
$ws.Range("B42").Value = "removeDeadBullets"

# --- The pre-existing asteroid subroutine rows (now at 43-47 after the
#     insert above) are reordered: spawnAsteroid, moveAsteroid,
#     drawAsteroids, removeAsteroid, valueToDirection8.
$ws.Range("B43").Value = "spawnAsteroid"
$ws.Range("B44").Value = "moveAsteroid"
$ws.Range("B45").Value = "drawAsteroids"
$ws.Range("B46").Value = "removeAsteroid"
$ws.Range("B47").Value = "valueToDirection8"

# --- Update the sheet view: scroll so row 17 is at the top, and select
#     B42 (the last-typed new subroutine name) as the active cell.
$ws.Activate()
$ws.Range("B42").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
